$d = $word.ActiveDocument

# Locate the paragraph that contains "Error message" (the last paragraph of the
# body, just before the final section properties).
$errPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Error message*") {
        $errPara = $p
    }
}

$targetRange = $errPara.Range
# Extend the range to also cover the paragraph mark so the replacement XML
# (which re-supplies an equivalent "Error message" paragraph as its first
# paragraph) cleanly swaps in without leaving a stray empty paragraph behind.
$targetRange.MoveEnd(1, 1)

$xmlPackage = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:cx="http://schemas.microsoft.com/office/drawing/2014/chartex" xmlns:cx1="http://schemas.microsoft.com/office/drawing/2015/9/8/chartex" xmlns:cx2="http://schemas.microsoft.com/office/drawing/2015/10/21/chartex" xmlns:cx3="http://schemas.microsoft.com/office/drawing/2016/5/9/chartex" xmlns:cx4="http://schemas.microsoft.com/office/drawing/2016/5/10/chartex" xmlns:cx5="http://schemas.microsoft.com/office/drawing/2016/5/11/chartex" xmlns:cx6="http://schemas.microsoft.com/office/drawing/2016/5/12/chartex" xmlns:cx7="http://schemas.microsoft.com/office/drawing/2016/5/13/chartex" xmlns:cx8="http://schemas.microsoft.com/office/drawing/2016/5/14/chartex" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:aink="http://schemas.microsoft.com/office/drawing/2016/ink" xmlns:am3d="http://schemas.microsoft.com/office/drawing/2017/model3d" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:w16cex="http://schemas.microsoft.com/office/word/2018/wordml/cex" xmlns:w16cid="http://schemas.microsoft.com/office/word/2016/wordml/cid" xmlns:w16="http://schemas.microsoft.com/office/word/2018/wordml" xmlns:w16se="http://schemas.microsoft.com/office/word/2015/wordml/symex" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 w16se w16cid w16 w16cex wp14"><w:body><w:p><w:pPr><w:ind w:left="720"/></w:pPr><w:r><w:t xml:space="preserve">     Error message</w:t></w:r></w:p>
<w:p/>
<w:p/>
<w:p>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="2800"/>
      <w:tab w:val="left" w:pos="5680"/>
    </w:tabs>
  </w:pPr>
  <w:r>
  <w:rPr>
    <w:noProof/>
  </w:rPr>
  <mc:AlternateContent>
    <mc:Choice Requires="wps">
      <w:drawing>
        <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251665408" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="46CCA52E" wp14:editId="33E9C13C">
          <wp:simplePos x="0" y="0"/>
          <wp:positionH relativeFrom="column">
            <wp:posOffset>1176867</wp:posOffset>
          </wp:positionH>
          <wp:positionV relativeFrom="paragraph">
            <wp:posOffset>149647</wp:posOffset>
          </wp:positionV>
          <wp:extent cx="507788" cy="465667"/>
          <wp:effectExtent l="0" t="0" r="64135" b="48895"/>
          <wp:wrapNone/>
          <wp:docPr id="7" name="Straight Arrow Connector 7"/>
          <wp:cNvGraphicFramePr/>
          <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
            <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
              <wps:wsp>
                <wps:cNvCnPr/>
                <wps:spPr>
                  <a:xfrm>
                    <a:off x="0" y="0"/>
                    <a:ext cx="507788" cy="465667"/>
                  </a:xfrm>
                  <a:prstGeom prst="straightConnector1">
                    <a:avLst/>
                  </a:prstGeom>
                  <a:ln>
                    <a:tailEnd type="triangle"/>
                  </a:ln>
                </wps:spPr>
                <wps:style>
                  <a:lnRef idx="1">
                    <a:schemeClr val="accent1"/>
                  </a:lnRef>
                  <a:fillRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:fillRef>
                  <a:effectRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:effectRef>
                  <a:fontRef idx="minor">
                    <a:schemeClr val="tx1"/>
                  </a:fontRef>
                </wps:style>
                <wps:bodyPr/>
              </wps:wsp>
            </a:graphicData>
          </a:graphic>
          <wp14:sizeRelH relativeFrom="margin">
            <wp14:pctWidth>0</wp14:pctWidth>
          </wp14:sizeRelH>
          <wp14:sizeRelV relativeFrom="margin">
            <wp14:pctHeight>0</wp14:pctHeight>
          </wp14:sizeRelV>
        </wp:anchor>
      </w:drawing>
    </mc:Choice>
    <mc:Fallback>
      <w:pict>
        <v:shapetype w14:anchorId="2FBC089A" id="_x0000_t32" coordsize="21600,21600" o:spt="32" o:oned="t" path="m,l21600,21600e" filled="f">
          <v:path arrowok="t" fillok="f" o:connecttype="none"/>
          <o:lock v:ext="edit" shapetype="t"/>
        </v:shapetype>
        <v:shape id="Straight Arrow Connector 7" o:spid="_x0000_s1026" type="#_x0000_t32" style="position:absolute;margin-left:92.65pt;margin-top:11.8pt;width:40pt;height:36.65pt;z-index:251665408;visibility:visible;mso-wrap-style:square;mso-width-percent:0;mso-height-percent:0;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text;mso-width-percent:0;mso-height-percent:0;mso-width-relative:margin;mso-height-relative:margin" strokecolor="#4472c4 [3204]" strokeweight=".5pt">
          <v:stroke endarrow="block" joinstyle="miter"/>
        </v:shape>
      </w:pict>
    </mc:Fallback>
  </mc:AlternateContent>
</w:r>
<w:r>
  <w:rPr>
    <w:noProof/>
  </w:rPr>
  <mc:AlternateContent>
    <mc:Choice Requires="wps">
      <w:drawing>
        <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251664384" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="29A6E760" wp14:editId="030719C8">
          <wp:simplePos x="0" y="0"/>
          <wp:positionH relativeFrom="column">
            <wp:posOffset>2810933</wp:posOffset>
          </wp:positionH>
          <wp:positionV relativeFrom="paragraph">
            <wp:posOffset>79587</wp:posOffset>
          </wp:positionV>
          <wp:extent cx="618067" cy="19261"/>
          <wp:effectExtent l="0" t="57150" r="10795" b="95250"/>
          <wp:wrapNone/>
          <wp:docPr id="6" name="Straight Arrow Connector 6"/>
          <wp:cNvGraphicFramePr/>
          <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
            <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
              <wps:wsp>
                <wps:cNvCnPr/>
                <wps:spPr>
                  <a:xfrm>
                    <a:off x="0" y="0"/>
                    <a:ext cx="618067" cy="19261"/>
                  </a:xfrm>
                  <a:prstGeom prst="straightConnector1">
                    <a:avLst/>
                  </a:prstGeom>
                  <a:ln>
                    <a:tailEnd type="triangle"/>
                  </a:ln>
                </wps:spPr>
                <wps:style>
                  <a:lnRef idx="1">
                    <a:schemeClr val="accent1"/>
                  </a:lnRef>
                  <a:fillRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:fillRef>
                  <a:effectRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:effectRef>
                  <a:fontRef idx="minor">
                    <a:schemeClr val="tx1"/>
                  </a:fontRef>
                </wps:style>
                <wps:bodyPr/>
              </wps:wsp>
            </a:graphicData>
          </a:graphic>
        </wp:anchor>
      </w:drawing>
    </mc:Choice>
    <mc:Fallback>
      <w:pict>
        <v:shape w14:anchorId="56867328" id="Straight Arrow Connector 6" o:spid="_x0000_s1026" type="#_x0000_t32" style="position:absolute;margin-left:221.35pt;margin-top:6.25pt;width:48.65pt;height:1.5pt;z-index:251664384;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text" strokecolor="#4472c4 [3204]" strokeweight=".5pt">
          <v:stroke endarrow="block" joinstyle="miter"/>
        </v:shape>
      </w:pict>
    </mc:Fallback>
  </mc:AlternateContent>
</w:r>
<w:r>
  <w:rPr>
    <w:noProof/>
  </w:rPr>
  <mc:AlternateContent>
    <mc:Choice Requires="wps">
      <w:drawing>
        <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251663360" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="4274C135" wp14:editId="1424E6E6">
          <wp:simplePos x="0" y="0"/>
          <wp:positionH relativeFrom="column">
            <wp:posOffset>1168400</wp:posOffset>
          </wp:positionH>
          <wp:positionV relativeFrom="paragraph">
            <wp:posOffset>73448</wp:posOffset>
          </wp:positionV>
          <wp:extent cx="516467" cy="8467"/>
          <wp:effectExtent l="0" t="57150" r="36195" b="86995"/>
          <wp:wrapNone/>
          <wp:docPr id="5" name="Straight Arrow Connector 5"/>
          <wp:cNvGraphicFramePr/>
          <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
            <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
              <wps:wsp>
                <wps:cNvCnPr/>
                <wps:spPr>
                  <a:xfrm>
                    <a:off x="0" y="0"/>
                    <a:ext cx="516467" cy="8467"/>
                  </a:xfrm>
                  <a:prstGeom prst="straightConnector1">
                    <a:avLst/>
                  </a:prstGeom>
                  <a:ln>
                    <a:tailEnd type="triangle"/>
                  </a:ln>
                </wps:spPr>
                <wps:style>
                  <a:lnRef idx="1">
                    <a:schemeClr val="accent1"/>
                  </a:lnRef>
                  <a:fillRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:fillRef>
                  <a:effectRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:effectRef>
                  <a:fontRef idx="minor">
                    <a:schemeClr val="tx1"/>
                  </a:fontRef>
                </wps:style>
                <wps:bodyPr/>
              </wps:wsp>
            </a:graphicData>
          </a:graphic>
        </wp:anchor>
      </w:drawing>
    </mc:Choice>
    <mc:Fallback>
      <w:pict>
        <v:shape w14:anchorId="6BF0A434" id="Straight Arrow Connector 5" o:spid="_x0000_s1026" type="#_x0000_t32" style="position:absolute;margin-left:92pt;margin-top:5.8pt;width:40.65pt;height:.65pt;z-index:251663360;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text" strokecolor="#4472c4 [3204]" strokeweight=".5pt">
          <v:stroke endarrow="block" joinstyle="miter"/>
        </v:shape>
      </w:pict>
    </mc:Fallback>
  </mc:AlternateContent>
</w:r>
<w:r>
  <w:rPr>
    <w:noProof/>
  </w:rPr>
  <w:t>Click</w:t>
</w:r>
<w:r>
  <w:t xml:space="preserve"> &#8220;Delete&#8221; Button </w:t>
</w:r>
<w:r>
  <w:tab/>
  <w:t xml:space="preserve">click &#8220;Yes&#8221; Button </w:t>
</w:r>
<w:r>
  <w:tab/>
</w:r>
<w:r>
  <w:t>Habit</w:t>
</w:r>
<w:r>
  <w:t xml:space="preserve"> is Deleted</w:t>
</w:r>
</w:p>
<w:p/>
<w:p>
  <w:pPr>
    <w:tabs>
      <w:tab w:val="left" w:pos="2760"/>
      <w:tab w:val="left" w:pos="5547"/>
    </w:tabs>
  </w:pPr>
  <w:r>
  <w:rPr>
    <w:noProof/>
  </w:rPr>
  <mc:AlternateContent>
    <mc:Choice Requires="wps">
      <w:drawing>
        <wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251666432" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1" wp14:anchorId="66BACB76" wp14:editId="271C59C4">
          <wp:simplePos x="0" y="0"/>
          <wp:positionH relativeFrom="column">
            <wp:posOffset>2760133</wp:posOffset>
          </wp:positionH>
          <wp:positionV relativeFrom="paragraph">
            <wp:posOffset>86148</wp:posOffset>
          </wp:positionV>
          <wp:extent cx="668655" cy="8467"/>
          <wp:effectExtent l="0" t="76200" r="17145" b="86995"/>
          <wp:wrapNone/>
          <wp:docPr id="8" name="Straight Arrow Connector 8"/>
          <wp:cNvGraphicFramePr/>
          <a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main">
            <a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape">
              <wps:wsp>
                <wps:cNvCnPr/>
                <wps:spPr>
                  <a:xfrm flipV="1">
                    <a:off x="0" y="0"/>
                    <a:ext cx="668655" cy="8467"/>
                  </a:xfrm>
                  <a:prstGeom prst="straightConnector1">
                    <a:avLst/>
                  </a:prstGeom>
                  <a:ln>
                    <a:tailEnd type="triangle"/>
                  </a:ln>
                </wps:spPr>
                <wps:style>
                  <a:lnRef idx="1">
                    <a:schemeClr val="accent1"/>
                  </a:lnRef>
                  <a:fillRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:fillRef>
                  <a:effectRef idx="0">
                    <a:schemeClr val="accent1"/>
                  </a:effectRef>
                  <a:fontRef idx="minor">
                    <a:schemeClr val="tx1"/>
                  </a:fontRef>
                </wps:style>
                <wps:bodyPr/>
              </wps:wsp>
            </a:graphicData>
          </a:graphic>
        </wp:anchor>
      </w:drawing>
    </mc:Choice>
    <mc:Fallback>
      <w:pict>
        <v:shape w14:anchorId="1E25E46E" id="Straight Arrow Connector 8" o:spid="_x0000_s1026" type="#_x0000_t32" style="position:absolute;margin-left:217.35pt;margin-top:6.8pt;width:52.65pt;height:.65pt;flip:y;z-index:251666432;visibility:visible;mso-wrap-style:square;mso-wrap-distance-left:9pt;mso-wrap-distance-top:0;mso-wrap-distance-right:9pt;mso-wrap-distance-bottom:0;mso-position-horizontal:absolute;mso-position-horizontal-relative:text;mso-position-vertical:absolute;mso-position-vertical-relative:text" strokecolor="#4472c4 [3204]" strokeweight=".5pt">
          <v:stroke endarrow="block" joinstyle="miter"/>
        </v:shape>
      </w:pict>
    </mc:Fallback>
  </mc:AlternateContent>
</w:r>
  <w:r>
  <w:tab/>
  <w:t>Click &#8220;No&#8221; Button</w:t>
</w:r>
<w:r>
  <w:tab/>
  <w:t>Page remains unchanged</w:t>
</w:r>
</w:p>
<w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($xmlPackage)

Write-Output "Inserted. Paragraph count now:"
Write-Output $d.Paragraphs.Count
